$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date value (serial 45182 = 2023-09-13) for every
# data row (2..472). The edit bumps that date forward by two days to serial
# 45184 (2023-09-15) across the whole column.
$ws.Range("C2:C472").Value = 45184
